# study font size and participant update
#
# - Participant's distractor condition for rows 25-26 (column B) changes
#   from "sit" to "stand".
# - Rows 4, 8, 12 and 16 get the same highlight formatting already used on
#   rows 2, 3, 6, 7, 10, 11, 14 and 15 (their order-block siblings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the distractor condition text.
$ws.Range("B25").Value = "stand"
$ws.Range("B26").Value = "stand"

# Match the highlight fill already applied to the alternating "block" rows.
$highlightColor = $ws.Range("A2:E2").Interior.Color
$ws.Range("A4:E4").Interior.Color = $highlightColor
$ws.Range("A8:E8").Interior.Color = $highlightColor
$ws.Range("A12:E12").Interior.Color = $highlightColor
$ws.Range("A16:E16").Interior.Color = $highlightColor

# Keep the last active selection / view in sync with the saved workbook.
$ws.Range("E23").Select()
